$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nemotecnico"
$ws.Range("B1").Value = "monto"
$ws.Range("C1").Value = "cantidad"
$ws.Range("D1").Value = "precio"
